$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 12.77268770170901
$ws.Cells.Item(2, 3).Value = 8.008463659254954
$ws.Cells.Item(2, 4).Value = 9.467861914974392
$ws.Cells.Item(2, 5).Value = 13.62946272699476
$ws.Cells.Item(2, 6).Value = 30.70433110987688
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 9).Value = 20.51028001455473
$ws.Cells.Item(2, 10).Value = 9.899118514847508
$ws.Cells.Item(2, 13).Value = 16.36564763522152
$ws.Cells.Item(2, 14).Value = 17.44350750119754
$ws.Cells.Item(2, 15).Value = 22.85327421657271
$ws.Cells.Item(3, 2).Value = 12.28844379339649
$ws.Cells.Item(3, 3).Value = 7.586567973294988
$ws.Cells.Item(3, 4).Value = 9.457719728724125
$ws.Cells.Item(3, 5).Value = 13.64549136321961
$ws.Cells.Item(3, 6).Value = 30.72502909154314
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 9).Value = 20.60406842943532
$ws.Cells.Item(3, 10).Value = 9.92371257768349
$ws.Cells.Item(3, 13).Value = 16.22393802648412
$ws.Cells.Item(3, 14).Value = 17.49257171515311
$ws.Cells.Item(3, 15).Value = 22.89958564744176
$ws.Cells.Item(4, 2).Value = 11.98255873707429
$ws.Cells.Item(4, 3).Value = 7.315542951004689
$ws.Cells.Item(4, 4).Value = 9.452769911162248
$ws.Cells.Item(4, 5).Value = 13.65750522157907
$ws.Cells.Item(4, 6).Value = 30.74583083609345
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 9).Value = 20.66612236587333
$ws.Cells.Item(4, 10).Value = 9.939934116059248
$ws.Cells.Item(4, 13).Value = 16.13878318282101
$ws.Cells.Item(4, 14).Value = 17.52444956731057
$ws.Cells.Item(4, 15).Value = 22.93369522678506
$ws.Cells.Item(5, 2).Value = 11.85594305670649
$ws.Cells.Item(5, 3).Value = 7.202197511780051
$ws.Cells.Item(5, 4).Value = 9.451075920913773
$ws.Cells.Item(5, 5).Value = 13.66294735658087
$ws.Cells.Item(5, 6).Value = 30.75634032176895
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 9).Value = 20.69253136154929
$ws.Cells.Item(5, 10).Value = 9.946826644733601
$ws.Cells.Item(5, 13).Value = 16.10457990110735
$ws.Cells.Item(5, 14).Value = 17.53788156545411
$ws.Cells.Item(5, 15).Value = 22.94901834871412
$ws.Cells.Item(6, 2).Value = 11.83480578874728
$ws.Cells.Item(6, 3).Value = 7.183205080059535
$ws.Cells.Item(6, 4).Value = 9.450814200735572
$ws.Cells.Item(6, 5).Value = 13.66388402241124
$ws.Cells.Item(6, 6).Value = 30.75820808461308
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 9).Value = 20.69698423763965
$ws.Cells.Item(6, 10).Value = 9.947988192524026
$ws.Cells.Item(6, 13).Value = 16.0989314525613
$ws.Cells.Item(6, 14).Value = 17.54013863172743
$ws.Cells.Item(6, 15).Value = 22.95164858460283
$ws.Cells.Item(7, 2).Value = 11.98085884854896
$ws.Cells.Item(7, 3).Value = 7.314025916031193
$ws.Cells.Item(7, 4).Value = 9.452745754855533
$ws.Cells.Item(7, 5).Value = 13.65757640365444
$ws.Cells.Item(7, 6).Value = 30.74596434485041
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 9).Value = 20.66647398832966
$ws.Cells.Item(7, 10).Value = 9.940025928398704
$ws.Cells.Item(7, 13).Value = 16.13831984809531
$ws.Cells.Item(7, 14).Value = 17.5246289270566
$ws.Cells.Item(7, 15).Value = 22.93389612243926
$ws.Cells.Item(8, 2).Value = 12.60760238019068
$ws.Cells.Item(8, 3).Value = 7.865541220452639
$ws.Cells.Item(8, 4).Value = 9.464100896948652
$ws.Cells.Item(8, 5).Value = 13.63453869026706
$ws.Cells.Item(8, 6).Value = 30.70978688691693
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 9).Value = 20.54168990985169
$ws.Cells.Item(8, 10).Value = 9.90736613914504
$ws.Cells.Item(8, 13).Value = 16.31642018386831
$ws.Cells.Item(8, 14).Value = 17.46006170892692
$ws.Cells.Item(8, 15).Value = 22.86806266419969
$ws.Cells.Item(9, 2).Value = 13.76153456963661
$ws.Cells.Item(9, 3).Value = 8.848056819868189
$ws.Cells.Item(9, 4).Value = 9.496417740375152
$ws.Cells.Item(9, 5).Value = 13.60658677707742
$ws.Cells.Item(9, 6).Value = 30.7031228021819
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 9).Value = 20.33251997102028
$ws.Cells.Item(9, 10).Value = 9.852198443268408
$ws.Cells.Item(9, 13).Value = 16.67897920019148
$ws.Cells.Item(9, 14).Value = 17.34730797661248
$ws.Cells.Item(9, 15).Value = 22.78413090059655
$ws.Cells.Item(10, 2).Value = 14.55518933605272
$ws.Cells.Item(10, 3).Value = 9.505375168991524
$ws.Cells.Item(10, 4).Value = 9.52616232920936
$ws.Cells.Item(10, 5).Value = 13.59653245585856
$ws.Cells.Item(10, 6).Value = 30.73743428077262
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 9).Value = 20.20061897495418
$ws.Cells.Item(10, 10).Value = 9.817059750133648
$ws.Cells.Item(10, 13).Value = 16.95152496174481
$ws.Cells.Item(10, 14).Value = 17.27286232805357
$ws.Cells.Item(10, 15).Value = 22.75017517359004
$ws.Cells.Item(11, 2).Value = 14.90306539094041
$ws.Cells.Item(11, 3).Value = 9.789770093215983
$ws.Cells.Item(11, 4).Value = 9.540967197720123
$ws.Cells.Item(11, 5).Value = 13.59422815910752
$ws.Cells.Item(11, 6).Value = 30.76153681830911
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 9).Value = 20.14536765739785
$ws.Cells.Item(11, 10).Value = 9.802241258767889
$ws.Cells.Item(11, 13).Value = 17.07643441589398
$ws.Cells.Item(11, 14).Value = 17.24080618717749
$ws.Cells.Item(11, 15).Value = 22.74077200883913
$ws.Cells.Item(12, 2).Value = 15.03280613820741
$ws.Cells.Item(12, 3).Value = 9.895321555292577
$ws.Cells.Item(12, 4).Value = 9.546753788058806
$ws.Cells.Item(12, 5).Value = 13.59368114416963
$ws.Cells.Item(12, 6).Value = 30.77188140559823
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 9).Value = 20.12513085567826
$ws.Cells.Item(12, 10).Value = 9.796797293030277
$ws.Cells.Item(12, 13).Value = 17.12383231464363
$ws.Cells.Item(12, 14).Value = 17.22892673254154
$ws.Cells.Item(12, 15).Value = 22.73808168012653
$ws.Cells.Item(13, 2).Value = 15.00495420349512
$ws.Cells.Item(13, 3).Value = 9.872685008286599
$ws.Cells.Item(13, 4).Value = 9.545499570534224
$ws.Cells.Item(13, 5).Value = 13.59378448745613
$ws.Cells.Item(13, 6).Value = 30.76959944239109
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 9).Value = 20.12945867793541
$ws.Cells.Item(13, 10).Value = 9.797962303637965
$ws.Cells.Item(13, 13).Value = 17.11362056710825
$ws.Cells.Item(13, 14).Value = 17.23147365481148
$ws.Cells.Item(13, 15).Value = 22.73862236064704
$ws.Cells.Item(14, 2).Value = 14.91377969376563
$ws.Cells.Item(14, 3).Value = 9.798497019048236
$ws.Cells.Item(14, 4).Value = 9.541439667388305
$ws.Cells.Item(14, 5).Value = 13.59417663581331
$ws.Cells.Item(14, 6).Value = 30.76236350079861
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 9).Value = 20.14368900740068
$ws.Cells.Item(14, 10).Value = 9.801790025303532
$ws.Cells.Item(14, 13).Value = 17.08033212259415
$ws.Cells.Item(14, 14).Value = 17.23982366026928
$ws.Cells.Item(14, 15).Value = 22.74053322273149
$ws.Cells.Item(15, 2).Value = 14.85767047947033
$ws.Cells.Item(15, 3).Value = 9.752774640461611
$ws.Cells.Item(15, 4).Value = 9.538976252822417
$ws.Cells.Item(15, 5).Value = 13.5944592118771
$ws.Cells.Item(15, 6).Value = 30.75808969033577
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 9).Value = 20.15249486367326
$ws.Cells.Item(15, 10).Value = 9.804156418691948
$ws.Cells.Item(15, 13).Value = 17.05995361050242
$ws.Cells.Item(15, 14).Value = 17.24497205513272
$ws.Cells.Item(15, 15).Value = 22.74181707109285
$ws.Cells.Item(16, 2).Value = 14.53218113448667
$ws.Cells.Item(16, 3).Value = 9.486491436027064
$ws.Cells.Item(16, 4).Value = 9.525220188020237
$ws.Cells.Item(16, 5).Value = 13.59672865137274
$ws.Cells.Item(16, 6).Value = 30.73602973628604
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 9).Value = 20.20432559348502
$ws.Cells.Item(16, 10).Value = 9.8180516218045
$ws.Cells.Item(16, 13).Value = 16.94337741684661
$ws.Cells.Item(16, 14).Value = 17.27499362468003
$ws.Cells.Item(16, 15).Value = 22.75091143466536
$ws.Cells.Item(17, 2).Value = 14.32905761750765
$ws.Cells.Item(17, 3).Value = 9.319358942802605
$ws.Cells.Item(17, 4).Value = 9.517105567754994
$ws.Cells.Item(17, 5).Value = 13.59870163475667
$ws.Cells.Item(17, 6).Value = 30.72466979058216
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 9).Value = 20.23734068947727
$ws.Cells.Item(17, 10).Value = 9.826874422917783
$ws.Cells.Item(17, 13).Value = 16.87207354727384
$ws.Cells.Item(17, 14).Value = 17.29387385010854
$ws.Cells.Item(17, 15).Value = 22.75803957671031
$ws.Cells.Item(18, 2).Value = 14.21099256471295
$ws.Cells.Item(18, 3).Value = 9.221854218464101
$ws.Cells.Item(18, 4).Value = 9.512558341149063
$ws.Cells.Item(18, 5).Value = 13.6000501100473
$ws.Cells.Item(18, 6).Value = 30.71893566975858
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 9).Value = 20.25677714683837
$ws.Cells.Item(18, 10).Value = 9.832058853997692
$ws.Cells.Item(18, 13).Value = 16.83115084836165
$ws.Cells.Item(18, 14).Value = 17.30490362116299
$ws.Cells.Item(18, 15).Value = 22.76270828476732
$ws.Cells.Item(19, 2).Value = 14.17080912926523
$ws.Cells.Item(19, 3).Value = 9.188606075284948
$ws.Cells.Item(19, 4).Value = 9.511039438762475
$ws.Cells.Item(19, 5).Value = 13.60054339960255
$ws.Cells.Item(19, 6).Value = 30.71713166465084
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 9).Value = 20.26343470332307
$ws.Cells.Item(19, 10).Value = 9.833833076862582
$ws.Cells.Item(19, 13).Value = 16.8173115764117
$ws.Cells.Item(19, 14).Value = 17.3086673941889
$ws.Cells.Item(19, 15).Value = 22.76438666771564
$ws.Cells.Item(20, 2).Value = 14.35080891992592
$ws.Cells.Item(20, 3).Value = 9.33729305757282
$ws.Cells.Item(20, 4).Value = 9.517956974599384
$ws.Cells.Item(20, 5).Value = 13.59846949871734
$ws.Cells.Item(20, 6).Value = 30.72579632251294
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 9).Value = 20.23377989059188
$ws.Cells.Item(20, 10).Value = 9.825923859605473
$ws.Cells.Item(20, 13).Value = 16.87965497377803
$ws.Cells.Item(20, 14).Value = 17.29184639091335
$ws.Cells.Item(20, 15).Value = 22.75722189437786
$ws.Cells.Item(21, 2).Value = 14.940614657793
$ws.Cells.Item(21, 3).Value = 9.820346261101751
$ws.Cells.Item(21, 4).Value = 9.542627289559345
$ws.Cells.Item(21, 5).Value = 13.59405262322749
$ws.Cells.Item(21, 6).Value = 30.76445586801838
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 9).Value = 20.13949058591614
$ws.Cells.Item(21, 10).Value = 9.800661187281037
$ws.Cells.Item(21, 13).Value = 17.09010738760545
$ws.Cells.Item(21, 14).Value = 17.23736402371489
$ws.Cells.Item(21, 15).Value = 22.73994832406414
$ws.Cells.Item(22, 2).Value = 15.31443673450961
$ws.Cells.Item(22, 3).Value = 10.12354696288454
$ws.Cells.Item(22, 4).Value = 9.559800295244006
$ws.Cells.Item(22, 5).Value = 13.59306317044494
$ws.Cells.Item(22, 6).Value = 30.79681578889584
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 9).Value = 20.08186486465046
$ws.Cells.Item(22, 10).Value = 9.785126638568615
$ws.Cells.Item(22, 13).Value = 17.22820141960693
$ws.Cells.Item(22, 14).Value = 17.20326896075469
$ws.Cells.Item(22, 15).Value = 22.73373293858851
$ws.Cells.Item(23, 2).Value = 15.11601637394319
$ws.Cells.Item(23, 3).Value = 9.962878185970808
$ws.Cells.Item(23, 4).Value = 9.550539697787192
$ws.Cells.Item(23, 5).Value = 13.59341795348191
$ws.Cells.Item(23, 6).Value = 30.77889725296063
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 9).Value = 20.11225417656545
$ws.Cells.Item(23, 10).Value = 9.793328485767045
$ws.Cells.Item(23, 13).Value = 17.15445941121227
$ws.Cells.Item(23, 14).Value = 17.22132798713882
$ws.Cells.Item(23, 15).Value = 22.73658562189202
$ws.Cells.Item(24, 2).Value = 14.34097915596238
$ws.Cells.Item(24, 3).Value = 9.329189458028718
$ws.Cells.Item(24, 4).Value = 9.517571685873303
$ws.Cells.Item(24, 5).Value = 13.59857378021688
$ws.Cells.Item(24, 6).Value = 30.72528453494911
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 9).Value = 20.2353883086439
$ws.Cells.Item(24, 10).Value = 9.826353260123208
$ws.Cells.Item(24, 13).Value = 16.8762271874303
$ws.Cells.Item(24, 14).Value = 17.29276245916976
$ws.Cells.Item(24, 15).Value = 22.75758979113048
$ws.Cells.Item(25, 2).Value = 13.45833705309719
$ws.Cells.Item(25, 3).Value = 8.593344972594998
$ws.Cells.Item(25, 4).Value = 9.486611571696974
$ws.Cells.Item(25, 5).Value = 13.61230591894635
$ws.Cells.Item(25, 6).Value = 30.69803711493694
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 9).Value = 20.3852903144539
$ws.Cells.Item(25, 10).Value = 9.866174316482027
$ws.Cells.Item(25, 13).Value = 16.57967132168928
$ws.Cells.Item(25, 14).Value = 17.37633269123383
$ws.Cells.Item(25, 15).Value = 22.80198136951048
